$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D (quarterly roll-forward: shifts old D:K data to F:M)
$ws.Range("D:E").Insert()

# Copy number formats / fonts / alignment from the (now shifted) F:G columns into the new D:E columns
$ws.Range("F:G").Copy()
$ws.Range("D:E").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns (D, E) with the latest reported figures
$ws.Range("D7").Value = 43465; $ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 114800; $ws.Range("E8").Value = 104800
$ws.Range("D9").Value = 49400; $ws.Range("E9").Value = 39200
$ws.Range("D10").Value = 65400; $ws.Range("E10").Value = 65600
$ws.Range("D12").Value = 21300; $ws.Range("E12").Value = 22000
$ws.Range("D13").Value = 0; $ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0; $ws.Range("E14").Value = 0
$ws.Range("D15").Value = "NA"; $ws.Range("E15").Value = "NA"
$ws.Range("D17").Value = 112800; $ws.Range("E17").Value = 100900
$ws.Range("D18").Value = 2000; $ws.Range("E18").Value = 3900
$ws.Range("D20").Value = 1000; $ws.Range("E20").Value = 1300
$ws.Range("D21").Value = 5400; $ws.Range("E21").Value = 8300
$ws.Range("D22").Value = 0; $ws.Range("E22").Value = 0
$ws.Range("D23").Value = 3000; $ws.Range("E23").Value = 5200
$ws.Range("D24").Value = 600; $ws.Range("E24").Value = -500
$ws.Range("D25").Value = 0; $ws.Range("E25").Value = 0
$ws.Range("D26").Value = 2400; $ws.Range("E26").Value = 5700
$ws.Range("D27").Value = 2400; $ws.Range("E27").Value = 5700
$ws.Range("D28").Value = 0; $ws.Range("E28").Value = 0
$ws.Range("D29").Value = -300; $ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0; $ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0; $ws.Range("E31").Value = 0
$ws.Range("D32").Value = -1000; $ws.Range("E32").Value = -1300
$ws.Range("D33").Value = 2100; $ws.Range("E33").Value = 5700
$ws.Range("D34").Value = 0; $ws.Range("E34").Value = 0
$ws.Range("D35").Value = 2100; $ws.Range("E35").Value = 5700
$ws.Range("D38").Value = 43465; $ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 349500; $ws.Range("E41").Value = 324400
$ws.Range("D42").Value = 0; $ws.Range("E42").Value = 500
$ws.Range("D43").Value = 144500; $ws.Range("E43").Value = 129800
$ws.Range("D44").Value = 33800; $ws.Range("E44").Value = 39200
$ws.Range("D45").Value = 30400; $ws.Range("E45").Value = 30500
$ws.Range("D46").Value = 558200; $ws.Range("E46").Value = 524400
$ws.Range("D47").Value = 40200; $ws.Range("E47").Value = 38200
$ws.Range("D48").Value = 37900; $ws.Range("E48").Value = 35600
$ws.Range("D49").Value = 40900; $ws.Range("E49").Value = 42000
$ws.Range("D50").Value = 0; $ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0; $ws.Range("E51").Value = 0
$ws.Range("D52").Value = 42300; $ws.Range("E52").Value = 41500
$ws.Range("D53").Value = 0; $ws.Range("E53").Value = 0
$ws.Range("D54").Value = 719500; $ws.Range("E54").Value = 681700
$ws.Range("D57").Value = 15200; $ws.Range("E57").Value = 9000
$ws.Range("D58").Value = 0; $ws.Range("E58").Value = 0
$ws.Range("D59").Value = 150800; $ws.Range("E59").Value = 132500
$ws.Range("D60").Value = 166000; $ws.Range("E60").Value = 141500
$ws.Range("D61").Value = 0; $ws.Range("E61").Value = 0
$ws.Range("D62").Value = 86200; $ws.Range("E62").Value = 80500
$ws.Range("D63").Value = 0; $ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0; $ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0; $ws.Range("E65").Value = 0
$ws.Range("D66").Value = 252200; $ws.Range("E66").Value = 222000
$ws.Range("D68").Value = 0; $ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0; $ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0; $ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0; $ws.Range("E71").Value = 0
$ws.Range("D72").Value = 171400; $ws.Range("E72").Value = 169300
$ws.Range("D73").Value = 0; $ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0; $ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0; $ws.Range("E75").Value = 0
$ws.Range("D76").Value = 467300; $ws.Range("E76").Value = 459700
$ws.Range("D77").Value = 0; $ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465; $ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 2100; $ws.Range("E81").Value = 5700
$ws.Range("D83").Value = 2400; $ws.Range("E83").Value = 3100
$ws.Range("D84").Value = 0; $ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0; $ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0; $ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0; $ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0; $ws.Range("E88").Value = 0
$ws.Range("D89").Value = 31200; $ws.Range("E89").Value = 16500
$ws.Range("D91").Value = -4300; $ws.Range("E91").Value = -2200
$ws.Range("D92").Value = 0; $ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0; $ws.Range("E93").Value = 0
$ws.Range("D94").Value = -3900; $ws.Range("E94").Value = 1200
$ws.Range("D96").Value = 0; $ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0; $ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0; $ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0; $ws.Range("E99").Value = 0
$ws.Range("D100").Value = -2800; $ws.Range("E100").Value = -1000
$ws.Range("D101").Value = -400; $ws.Range("E101").Value = 200
$ws.Range("D102").Value = 24200; $ws.Range("E102").Value = 16900

# A handful of prior-quarter figures were also restated in this update
$ws.Range("F91").Value = -3600
$ws.Range("H91").Value = -1300
$ws.Range("I91").Value = -3300
$ws.Range("J91").Value = -3400
$ws.Range("H100").Value = -2200
$ws.Range("I100").Value = -200
